# Generate Report for Handoff
# Updates the status/handoff info for file "b.md" now that it has been
# handed off (b.63290e5768f688058c7b37413b0a5c26c308f864 content id).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
# Row 3 corresponds to b.md
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-23 16:37:02"

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
# Row 3 corresponds to b.md
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-23 16:36:58"

# fix up the hyperlink display text on D3 (target address is unchanged)
$zhcnTarget = $null
foreach ($hl in $zhcn.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$D`$3") {
        $zhcnTarget = $hl
    }
}
if ($zhcnTarget -ne $null) {
    $zhcnAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2b17c8aee8696b60a656652f85665866115e25ce/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
    $zhcnTarget.Delete()
    $zhcn.Hyperlinks.Add($zhcn.Range("D3"), $zhcnAddr, "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf") | Out-Null
    # keep the cell looking like the other hyperlink cells (blue underline)
    $zhcn.Range("D3").Font.Name = "Calibri"
    $zhcn.Range("D3").Font.Size = 11
    $zhcn.Range("D3").Font.Underline = 2
    $zhcn.Range("D3").Font.Color = 15570276
}

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
# Row 3 corresponds to b.md
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-23 16:37:02"

# fix up the hyperlink display text on D3 (target address is unchanged)
$dedeTarget = $null
foreach ($hl in $dede.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$D`$3") {
        $dedeTarget = $hl
    }
}
if ($dedeTarget -ne $null) {
    $dedeAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9a39f568a3b65f73167a51bcd513cbe73a22ec82/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
    $dedeTarget.Delete()
    $dede.Hyperlinks.Add($dede.Range("D3"), $dedeAddr, "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf") | Out-Null
    # keep the cell looking like the other hyperlink cells (blue underline)
    $dede.Range("D3").Font.Name = "Calibri"
    $dede.Range("D3").Font.Size = 11
    $dede.Range("D3").Font.Underline = 2
    $dede.Range("D3").Font.Color = 15570276
}
